# "create new slide 2"
# Insert a new slide at position 2 using the "Title and Content" auto
# layout (legacy PowerPoint layout id 2 == ppLayoutText), which maps to
# this deck's second slide layout ("Заголовок и объект" / slideLayout2.xml)
# -- a title placeholder plus a generic content placeholder (idx=1).
# Both placeholders are left with their default empty text, matching a
# freshly inserted, untouched slide.

$p = $ppt.ActivePresentation
$ppLayoutText = 2
$s = $p.Slides.Add(2, $ppLayoutText)
